$wb = $excel.ActiveWorkbook

# --- ALC (sheet1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 154.7
$ws.Range("I6").Value = 149.66667
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 449.00001
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -337.00001
$ws.Range("N6").Value = -824
$ws.Range("H38").Value = 1275.6666
$ws.Range("I38").Value = 538.5
$ws.Range("K38").Value = 1615.5
$ws.Range("M38").Value = -1243.5
$ws.Range("H39").Value = 24
$ws.Range("I39").Value = 24
$ws.Range("K39").Value = 72
$ws.Range("M39").Value = 224
$ws.Range("H53").Value = 349
$ws.Range("I53").Value = 95
$ws.Range("K53").Value = 95
$ws.Range("M53").Value = 542
$ws.Range("H98").Value = 913.13635
$ws.Range("I98").Value = 904.5
$ws.Range("K98").Value = 904.5
$ws.Range("M98").Value = 593.5
$ws.Range("H122").Value = 913.13635
$ws.Range("I122").Value = 904.5
$ws.Range("K122").Value = 2713.5
$ws.Range("M122").Value = -263.5

# --- ARM (sheet2) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 327.625
$ws.Range("I4").Value = 283
$ws.Range("J4").Value = 354.4
$ws.Range("K4").Value = 283
$ws.Range("L4").Value = 354.4
$ws.Range("M4").Value = -167
$ws.Range("N4").Value = -586.4
$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = -474
$ws.Range("H61").Value = 1966.3077
$ws.Range("I61").Value = 1936.826
$ws.Range("K61").Value = 1936.826
$ws.Range("M61").Value = -1724.826
$ws.Range("H74").Value = 647.5333000000001
$ws.Range("I74").Value = 550.9286
$ws.Range("K74").Value = 550.9286
$ws.Range("M74").Value = 323.0714
$ws.Range("H77").Value = 647.5333000000001
$ws.Range("I77").Value = 550.9286
$ws.Range("K77").Value = 2754.643
$ws.Range("M77").Value = 1613.357
$ws.Range("H102").Value = 6947108
$ws.Range("I102").Value = 15875076
$ws.Range("J102").Value = 3133.3333
$ws.Range("K102").Value = 15875076
$ws.Range("L102").Value = 3133.3333
$ws.Range("M102").Value = -15873454
$ws.Range("N102").Value = -6377.3333
$ws.Range("H110").Value = 10103484
$ws.Range("I110").Value = 10103484
$ws.Range("K110").Value = 10103484
$ws.Range("M110").Value = -10101439
$ws.Range("H122").Value = 11321.667
$ws.Range("J122").Value = 20333.334
$ws.Range("L122").Value = 61000.00199999999
$ws.Range("N122").Value = -65900.00199999999
$ws.Range("H132").Value = 2355.5715
$ws.Range("I132").Value = 2152.1538
$ws.Range("K132").Value = 6456.4614
$ws.Range("M132").Value = -3926.4614
$ws.Range("H136").Value = 1966.3077
$ws.Range("I136").Value = 1936.826
$ws.Range("K136").Value = 5810.478
$ws.Range("M136").Value = -3260.478

# --- BSM (sheet3) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -480
$ws.Range("H20").Value = 3499.5
$ws.Range("I20").Value = 3499.5
$ws.Range("K20").Value = 3499.5
$ws.Range("M20").Value = -3252.5
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null
$ws.Range("H94").Value = 1753.7858
$ws.Range("I94").Value = 1109.5
$ws.Range("J94").Value = 3364.5
$ws.Range("K94").Value = 1109.5
$ws.Range("L94").Value = 3364.5
$ws.Range("M94").Value = -658.5
$ws.Range("N94").Value = -4266.5
$ws.Range("H105").Value = 4389117
$ws.Range("I105").Value = 8336285
$ws.Range("K105").Value = 8336285
$ws.Range("M105").Value = -8334538
$ws.Range("H107").Value = 3229.125
$ws.Range("I107").Value = 2445
$ws.Range("J107").Value = 4536
$ws.Range("K107").Value = 2445
$ws.Range("L107").Value = 4536
$ws.Range("M107").Value = -525
$ws.Range("N107").Value = -8376
$ws.Range("H134").Value = 1771.8235
$ws.Range("I134").Value = 1205.1428
$ws.Range("K134").Value = 3615.4284
$ws.Range("M134").Value = -1080.4284

# --- CRP (sheet4) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 242.25
$ws.Range("I7").Value = 243
$ws.Range("J7").Value = 240
$ws.Range("K7").Value = 243
$ws.Range("L7").Value = 240
$ws.Range("M7").Value = -130
$ws.Range("N7").Value = -466
$ws.Range("H15").Value = 10809.333
$ws.Range("I15").Value = 11710
$ws.Range("K15").Value = 11710
$ws.Range("M15").Value = -11540
$ws.Range("H99").Value = 13437.261
$ws.Range("I99").Value = 10204.6
$ws.Range("J99").Value = 15923.923
$ws.Range("K99").Value = 10204.6
$ws.Range("L99").Value = 15923.923
$ws.Range("M99").Value = -8706.6
$ws.Range("N99").Value = -18919.923
$ws.Range("H105").Value = 2978.7896
$ws.Range("I105").Value = 1290.3334
$ws.Range("K105").Value = 1290.3334
$ws.Range("M105").Value = 456.6666
$ws.Range("H126").Value = 13437.261
$ws.Range("I126").Value = 10204.6
$ws.Range("J126").Value = 15923.923
$ws.Range("K126").Value = 30613.8
$ws.Range("L126").Value = 47771.769
$ws.Range("M126").Value = -28143.8
$ws.Range("N126").Value = -52711.769

# --- GSM (sheet6) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 266.57144
$ws.Range("I2").Value = 49.875
$ws.Range("K2").Value = 49.875
$ws.Range("M2").Value = 63.125
$ws.Range("H55").Value = 7000
$ws.Range("J55").Value = 7000
$ws.Range("L55").Value = 7000
$ws.Range("N55").Value = -7654
$ws.Range("H113").Value = 79998.336
$ws.Range("J113").Value = 79998.336
$ws.Range("L113").Value = 79998.336
$ws.Range("N113").Value = -84338.336
$ws.Range("H126").Value = 4941
$ws.Range("J126").Value = 4949.5
$ws.Range("L126").Value = 14848.5
$ws.Range("N126").Value = -19788.5
$ws.Range("H132").Value = 2427
$ws.Range("I132").Value = 1890.5
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 5671.5
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -3141.5
$ws.Range("N132").Value = -15560

# --- LTW (sheet7) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = $null
$ws.Range("H93").Value = 1206.8334
$ws.Range("I93").Value = 1108
$ws.Range("K93").Value = 1108
$ws.Range("M93").Value = 140

# --- WVR (sheet8) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3314.3845
$ws.Range("I132").Value = 2414.1428
$ws.Range("J132").Value = 4364.6665
$ws.Range("K132").Value = 7242.428400000001
$ws.Range("L132").Value = 13093.9995
$ws.Range("M132").Value = -4712.428400000001
$ws.Range("N132").Value = -18153.9995
$ws.Range("H135").Value = 78995
$ws.Range("J135").Value = 78995
$ws.Range("L135").Value = 78995
$ws.Range("N135").Value = -89135
$ws.Range("H136").Value = 1316.5
$ws.Range("I136").Value = 842.94116
$ws.Range("K136").Value = 2528.82348
$ws.Range("M136").Value = 21.17651999999998
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null
